# Auto-generated Excel COM-interop script applying market-price refresh edits
# to the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 225
$ws.Range("I6").Value = 179
$ws.Range("J6").Value = 261.8
$ws.Range("K6").Value = 537
$ws.Range("L6").Value = 785.4000000000001
$ws.Range("M6").Value = -425

$ws.Range("H40").Value = 3623.5
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 3623.5
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 3623.5
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -3973.5

$ws.Range("H43").Value = 589364.3
$ws.Range("I43").Value = 2975
$ws.Range("J43").Value = 823920
$ws.Range("K43").Value = 2975
$ws.Range("L43").Value = 823920
$ws.Range("M43").Value = -2906
$ws.Range("N43").Value = -824058

$ws.Range("H55").Value = 653.5714
$ws.Range("I55").Value = 745.8333
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 745.8333
$ws.Range("L55").Value = 100
$ws.Range("M55").Value = -531.8333
$ws.Range("N55").Value = -528

$ws.Range("H98").Value = 2893.4546
$ws.Range("I98").Value = 2893.4546
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 2893.4546
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = -1395.4546
$ws.Range("N98").ClearContents()

$ws.Range("H122").Value = 2893.4546
$ws.Range("I122").Value = 2893.4546
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8680.363799999999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6230.363799999999
$ws.Range("N122").ClearContents()

$ws.Range("H131").Value = 3334.0908
$ws.Range("I131").Value = 1065.5264
$ws.Range("J131").Value = 17701.666
$ws.Range("K131").Value = 3196.5792
$ws.Range("L131").Value = 53104.99800000001
$ws.Range("M131").Value = 1843.4208
$ws.Range("N131").Value = -63184.99800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 31939.912
$ws.Range("I74").Value = 45132.957
$ws.Range("J74").Value = 4354.4546
$ws.Range("K74").Value = 45132.957
$ws.Range("L74").Value = 4354.4546
$ws.Range("M74").Value = -44258.957

$ws.Range("H77").Value = 31939.912
$ws.Range("I77").Value = 45132.957
$ws.Range("J77").Value = 4354.4546
$ws.Range("K77").Value = 225664.785
$ws.Range("L77").Value = 21772.273
$ws.Range("M77").Value = -221296.785

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 2450
$ws.Range("I30").Value = 2450
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2450
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -2325

$ws.Range("H105").Value = 3433.5
$ws.Range("I105").Value = 2608.5217
$ws.Range("J105").Value = 4432.1577
$ws.Range("K105").Value = 2608.5217
$ws.Range("L105").Value = 4432.1577
$ws.Range("M105").Value = -861.5216999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 249.07143
$ws.Range("I7").Value = 123.333336
$ws.Range("J7").Value = 343.375
$ws.Range("K7").Value = 123.333336
$ws.Range("L7").Value = 343.375
$ws.Range("M7").Value = -10.333336
$ws.Range("N7").Value = -569.375

$ws.Range("H58").Value = 5792.8887
$ws.Range("I58").Value = 2024.762
$ws.Range("J58").Value = 9090
$ws.Range("K58").Value = 2024.762
$ws.Range("L58").Value = 9090
$ws.Range("M58").Value = -1821.762

$ws.Range("H98").Value = 55001
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 55001
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 55001
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -59493

$ws.Range("H136").Value = 5792.8887
$ws.Range("I136").Value = 2024.762
$ws.Range("J136").Value = 9090
$ws.Range("K136").Value = 6074.286
$ws.Range("L136").Value = 27270
$ws.Range("M136").Value = -3524.286

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H48").Value = 8328.166999999999
$ws.Range("I48").Value = 0
$ws.Range("J48").Value = 8328.166999999999
$ws.Range("K48").Value = 0
$ws.Range("L48").Value = 24984.501
$ws.Range("N48").Value = -25484.501

$ws.Range("H134").Value = 76250.14
$ws.Range("I134").Value = 81346.30499999999
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 244038.915
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -238968.915

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 66894.25
$ws.Range("I62").Value = 30077
$ws.Range("J62").Value = 79166.664
$ws.Range("K62").Value = 30077
$ws.Range("L62").Value = 79166.664
$ws.Range("M62").Value = -29391
$ws.Range("N62").Value = -80538.664

$ws.Range("H65").Value = 66894.25
$ws.Range("I65").Value = 30077
$ws.Range("J65").Value = 79166.664
$ws.Range("K65").Value = 90231
$ws.Range("L65").Value = 237499.992
$ws.Range("M65").Value = -86799
$ws.Range("N65").Value = -244363.992

$ws.Range("H102").Value = 4487.069
$ws.Range("I102").Value = 4295.636
$ws.Range("J102").Value = 5088.7144
$ws.Range("K102").Value = 4295.636
$ws.Range("L102").Value = 5088.7144
$ws.Range("M102").Value = -2673.636

$ws.Range("H132").Value = 7750.0835
$ws.Range("I132").Value = 3001.7144
$ws.Range("J132").Value = 14397.8
$ws.Range("K132").Value = 9005.143199999999
$ws.Range("L132").Value = 43193.39999999999
$ws.Range("M132").Value = -6475.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5910
$ws.Range("I7").Value = 4308.7144
$ws.Range("J7").Value = 8151.8
$ws.Range("K7").Value = 4308.7144
$ws.Range("L7").Value = 8151.8
$ws.Range("M7").Value = -4196.7144
$ws.Range("N7").Value = -8375.799999999999

$ws.Range("H20").Value = 2400000
$ws.Range("I20").Value = 2366666.8
$ws.Range("J20").Value = 2500000
$ws.Range("K20").Value = 2366666.8
$ws.Range("L20").Value = 2500000
$ws.Range("M20").Value = -2366440.8
$ws.Range("N20").Value = -2500452

$ws.Range("H40").Value = 2836.8
$ws.Range("I40").Value = 3128
$ws.Range("J40").Value = 2400
$ws.Range("K40").Value = 3128
$ws.Range("L40").Value = 2400
$ws.Range("M40").Value = -2992

$ws.Range("H68").Value = 4464.2383
$ws.Range("I68").Value = 2645
$ws.Range("J68").Value = 6118.091
$ws.Range("K68").Value = 2645
$ws.Range("L68").Value = 6118.091
$ws.Range("M68").Value = -1896
$ws.Range("N68").Value = -7616.091

$ws.Range("H71").Value = 4464.2383
$ws.Range("I71").Value = 2645
$ws.Range("J71").Value = 6118.091
$ws.Range("K71").Value = 13225
$ws.Range("L71").Value = 30590.455
$ws.Range("M71").Value = -9481
$ws.Range("N71").Value = -38078.455

$ws.Range("H122").Value = 4254.636
$ws.Range("I122").Value = 2692.6924
$ws.Range("J122").Value = 6510.778
$ws.Range("K122").Value = 8078.0772
$ws.Range("L122").Value = 19532.334
$ws.Range("M122").Value = -5628.0772
$ws.Range("N122").Value = -24432.334

$ws.Range("H126").Value = 5910
$ws.Range("I126").Value = 4308.7144
$ws.Range("J126").Value = 8151.8
$ws.Range("K126").Value = 12926.1432
$ws.Range("L126").Value = 24455.4
$ws.Range("M126").Value = -10456.1432
$ws.Range("N126").Value = -29395.4

$ws.Range("H132").Value = 7053.7803
$ws.Range("I132").Value = 3520.4
$ws.Range("J132").Value = 9092.27
$ws.Range("K132").Value = 10561.2
$ws.Range("L132").Value = 27276.81
$ws.Range("M132").Value = -8031.200000000001

$ws.Range("H136").Value = 17396.465
$ws.Range("I136").Value = 3513
$ws.Range("J136").Value = 22949.85
$ws.Range("K136").Value = 10539
$ws.Range("L136").Value = 68849.54999999999
$ws.Range("M136").Value = -7989

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H30").Value = 12500
$ws.Range("I30").Value = 10000
$ws.Range("J30").Value = 15000
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 15000
$ws.Range("M30").Value = -9893
$ws.Range("N30").Value = -15214

$ws.Range("H70").Value = 30599.166
$ws.Range("I70").Value = 20095
$ws.Range("J70").Value = 32700
$ws.Range("K70").Value = 20095
$ws.Range("L70").Value = 32700
$ws.Range("M70").Value = -19780
$ws.Range("N70").Value = -33330

$ws.Range("H73").Value = 30599.166
$ws.Range("I73").Value = 20095
$ws.Range("J73").Value = 32700
$ws.Range("K73").Value = 20095
$ws.Range("L73").Value = 32700
$ws.Range("M73").Value = -19003
$ws.Range("N73").Value = -34884

$ws.Range("H96").Value = 2999
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 2999
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 2999
$ws.Range("N96").Value = -5745

$ws.Range("H107").Value = 18519866
$ws.Range("I107").Value = 614
$ws.Range("J107").Value = 30304844
$ws.Range("K107").Value = 1842
$ws.Range("L107").Value = 90914532
$ws.Range("M107").Value = 78

$ws.Range("H132").Value = 6217.5835
$ws.Range("I132").Value = 9429.454
$ws.Range("J132").Value = 3499.8462
$ws.Range("K132").Value = 28288.362
$ws.Range("L132").Value = 10499.5386
$ws.Range("M132").Value = -25758.362
$ws.Range("N132").Value = -15559.5386

$ws.Range("H133").Value = 170975.2
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 170975.2
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 170975.2
$ws.Range("N133").Value = -181095.2

$ws.Range("H136").Value = 405629.9
$ws.Range("I136").Value = 1915.6666
$ws.Range("J136").Value = 632719.2
$ws.Range("K136").Value = 5746.9998
$ws.Range("L136").Value = 1898157.6
$ws.Range("M136").Value = -3196.9998
$ws.Range("N136").Value = -1903257.6
